$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 4727273
$ws.Range("C4").Value = 21384
$ws.Range("D4").Value = 2329402
$ws.Range("E4").Value = 2240701
$ws.Range("G4").Value = 423
$ws.Range("H4").Value = 157170

# --- Row 6: India ---
$ws.Range("B6").Value = 1749771
$ws.Range("C6").Value = 52717
$ws.Range("D6").Value = 1144277
$ws.Range("E6").Value = 568104
$ws.Range("G6").Value = 839
$ws.Range("H6").Value = 37390

# --- Row 11: Chile ---
$ws.Range("B11").Value = 357658
$ws.Range("C11").Value = 1991
$ws.Range("D11").Value = 330507
$ws.Range("E11").Value = 17618
$ws.Range("G11").Value = 76
$ws.Range("H11").Value = 9533

# --- Row 14: Reino Unido ---
$ws.Range("G14").Value = 74
$ws.Range("H14").Value = 46193

# --- Row 18: Italia ---
$ws.Range("B18").Value = 247832
$ws.Range("C18").Value = 295
$ws.Range("D18").Value = 200229
$ws.Range("E18").Value = 12457
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 35146

# --- Row 45: Singapur ---
$ws.Range("D45").Value = 46740
$ws.Range("E45").Value = 5745

# --- Row 75: Chequia ---
$ws.Range("B75").Value = 16640
$ws.Range("C75").Value = 232
$ws.Range("D75").Value = 11584
$ws.Range("E75").Value = 4674

# --- Rows 115/116: Zimbabue & Montenegro swap ranking order ---
# Row 115 becomes Montenegro (was Zimbabue), with new stats
$ws.Range("A115").Value = "Montenegro"
$ws.Range("B115").Value = 3198
$ws.Range("C115").Value = 86
$ws.Range("D115").Value = 1293
$ws.Range("E115").Value = 1855
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 50

# Row 116 becomes Zimbabue (was Montenegro), carrying the old Zimbabue stats
$ws.Range("A116").Value = "Zimbabue"
$ws.Range("B116").Value = 3169
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 1004
$ws.Range("E116").Value = 2098
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 67

# --- Row 120: Cuba ---
$ws.Range("B120").Value = 2633
$ws.Range("C120").Value = 25
$ws.Range("D120").Value = 2367
$ws.Range("E120").Value = 179

# --- Rows 158/159/160: San Marino, Reunion, Lesoto re-ranked ---
# Row 158 becomes Lesoto (was San Marino), with new stats
$ws.Range("A158").Value = "Lesoto"
$ws.Range("B158").Value = 702
$ws.Range("C158").Value = 98
$ws.Range("D158").Value = 171
$ws.Range("E158").Value = 517
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 14

# Row 159 becomes San Marino (was Reunion), carrying old San Marino stats
$ws.Range("A159").Value = "San Marino"
$ws.Range("B159").Value = 699
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 657
$ws.Range("E159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 42

# Row 160 becomes Reunion (was Lesoto), carrying old Reunion stats
$ws.Range("A160").Value = "Reunion"
$ws.Range("B160").Value = 660
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 592
$ws.Range("E160").Value = 64
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 4

# --- Timestamp update (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 18:06"
